$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.141.78"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.408.25"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "558.28"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").Value = "135.26"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "0.349"
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("D13").Value = "24.71"
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("D14").Value = "2.838.45"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "60.064.93"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "2.346.16"
$ws.Range("E17").Value = "  -3.29%  "
$ws.Range("D18").Value = "11.22"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("E19").Value = "  +2.71%  "
$ws.Range("D20").Value = "326.22"
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("D21").Value = "6.80"
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "64.77"
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -2.83%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "1.39"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").Value = "170.64"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("E32").Value = "  +5.25%  "
$ws.Range("D33").Value = "0.402"
$ws.Range("E33").Value = "  -2.15%  "
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.33"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").Value = "324.49"
$ws.Range("E39").Value = "  +3.52%  "
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").Value = "38.54"
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("D42").Value = "148.65"
$ws.Range("E42").Value = "  +6.83%  "
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("D44").Value = "0.0968"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").Value = "19.93"
$ws.Range("E45").Value = "  +2.22%  "
$ws.Range("D46").Value = "0.0516"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").Value = "0.576"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("D51").Value = "4.67"
$ws.Range("E51").Value = "  -0.82%  "
